$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.5049972534179688
$ws.Cells.Item(2, 5).Value = 109.0979712069511
$ws.Cells.Item(2, 6).Value = 0.004056021945133336
$ws.Cells.Item(2, 7).Value = 0.003174990572092158
$ws.Cells.Item(2, 8).Value = 0.00289868642178024
$ws.Cells.Item(2, 9).Value = 0.002752399711127126
$ws.Cells.Item(2, 10).Value = 0.002752399711127126
$ws.Cells.Item(2, 11).Value = 0.00263654202077507
$ws.Cells.Item(2, 12).Value = 0.002605533143348552
$ws.Cells.Item(2, 13).Value = 0.002515283366239146
$ws.Cells.Item(2, 14).Value = 0.002436642529629436
$ws.Cells.Item(2, 15).Value = 0.002299545719419554
$ws.Cells.Item(2, 16).Value = 0.002295255363988849
$ws.Cells.Item(2, 17).Value = 0.002295255363988849
$ws.Cells.Item(2, 18).Value = 0.002257737240391057
$ws.Cells.Item(2, 19).Value = 0.002223139801280851
$ws.Cells.Item(2, 20).Value = 0.002171425944798492
$ws.Cells.Item(2, 21).Value = 0.002171425944798492
$ws.Cells.Item(2, 22).Value = 0.002171425944798492
$ws.Cells.Item(2, 23).Value = 0.002146709201854029
$ws.Cells.Item(2, 24).Value = 0.002128931898614415
$ws.Cells.Item(2, 25).Value = 0.002126666105398657

$ws.Cells.Item(3, 3).Value = 0.3930010795593262
$ws.Cells.Item(3, 5).Value = 109.7886991023606
$ws.Cells.Item(3, 6).Value = 0.003879546177076297
$ws.Cells.Item(3, 7).Value = 0.003185387997732425
$ws.Cells.Item(3, 8).Value = 0.003185387997732425
$ws.Cells.Item(3, 9).Value = 0.003013914824927974
$ws.Cells.Item(3, 10).Value = 0.002878521169169585
$ws.Cells.Item(3, 11).Value = 0.002716485241637899
$ws.Cells.Item(3, 12).Value = 0.002567788796163209
$ws.Cells.Item(3, 13).Value = 0.002453647407975957
$ws.Cells.Item(3, 14).Value = 0.002453647407975957
$ws.Cells.Item(3, 15).Value = 0.002424460771586704
$ws.Cells.Item(3, 16).Value = 0.002365908120050786
$ws.Cells.Item(3, 17).Value = 0.002322238994612158
$ws.Cells.Item(3, 18).Value = 0.002306683638556468
$ws.Cells.Item(3, 19).Value = 0.002246189554539119
$ws.Cells.Item(3, 20).Value = 0.002238091579048934
$ws.Cells.Item(3, 21).Value = 0.002227069722937256
$ws.Cells.Item(3, 22).Value = 0.002187058597320293
$ws.Cells.Item(3, 23).Value = 0.002165339396188927
$ws.Cells.Item(3, 24).Value = 0.002160262283776259
$ws.Cells.Item(3, 25).Value = 0.002140130586790655

$ws.Cells.Item(4, 3).Value = 0.3709585666656494
$ws.Cells.Item(4, 5).Value = 107.8352910095036
$ws.Cells.Item(4, 6).Value = 0.00410297909976543
$ws.Cells.Item(4, 7).Value = 0.003194297563566305
$ws.Cells.Item(4, 8).Value = 0.003035415894435353
$ws.Cells.Item(4, 9).Value = 0.002981810686580881
$ws.Cells.Item(4, 10).Value = 0.002790133411989776
$ws.Cells.Item(4, 11).Value = 0.002788073239641737
$ws.Cells.Item(4, 12).Value = 0.002788073239641737
$ws.Cells.Item(4, 13).Value = 0.002547773169650989
$ws.Cells.Item(4, 14).Value = 0.002477368537889262
$ws.Cells.Item(4, 15).Value = 0.002477044301649194
$ws.Cells.Item(4, 16).Value = 0.002417902520391948
$ws.Cells.Item(4, 17).Value = 0.002302179009498591
$ws.Cells.Item(4, 18).Value = 0.002247334813035631
$ws.Cells.Item(4, 19).Value = 0.002191282436098923
$ws.Cells.Item(4, 20).Value = 0.002191282436098923
$ws.Cells.Item(4, 21).Value = 0.002191282436098923
$ws.Cells.Item(4, 22).Value = 0.002175117926101938
$ws.Cells.Item(4, 23).Value = 0.002149627195426023
$ws.Cells.Item(4, 24).Value = 0.002102052456325606
$ws.Cells.Item(4, 25).Value = 0.002102052456325606

$ws.Cells.Item(5, 3).Value = 0.3780362606048584
$ws.Cells.Item(5, 5).Value = 110.4816627382042
$ws.Cells.Item(5, 6).Value = 0.003917291961539131
$ws.Cells.Item(5, 7).Value = 0.003198531569195744
$ws.Cells.Item(5, 8).Value = 0.00304435214499877
$ws.Cells.Item(5, 9).Value = 0.00292316521416381
$ws.Cells.Item(5, 10).Value = 0.002766823730599446
$ws.Cells.Item(5, 11).Value = 0.002577130377100231
$ws.Cells.Item(5, 12).Value = 0.002577130377100231
$ws.Cells.Item(5, 13).Value = 0.002510442361909637
$ws.Cells.Item(5, 14).Value = 0.002510442361909637
$ws.Cells.Item(5, 15).Value = 0.002510442361909637
$ws.Cells.Item(5, 16).Value = 0.002421795447260709
$ws.Cells.Item(5, 17).Value = 0.002421795447260709
$ws.Cells.Item(5, 18).Value = 0.0023116567736373
$ws.Cells.Item(5, 19).Value = 0.002267157208067221
$ws.Cells.Item(5, 20).Value = 0.002211449358321728
$ws.Cells.Item(5, 21).Value = 0.002197362022840276
$ws.Cells.Item(5, 22).Value = 0.002190676528939964
$ws.Cells.Item(5, 23).Value = 0.002155051771358941
$ws.Cells.Item(5, 24).Value = 0.002155051771358941
$ws.Cells.Item(5, 25).Value = 0.002153638649867529

$ws.Cells.Item(6, 3).Value = 0.4140205383300781
$ws.Cells.Item(6, 5).Value = 106.4137722306659
$ws.Cells.Item(6, 6).Value = 0.003943007570197137
$ws.Cells.Item(6, 7).Value = 0.003178116674314109
$ws.Cells.Item(6, 8).Value = 0.002961732077310762
$ws.Cells.Item(6, 9).Value = 0.00273004666891374
$ws.Cells.Item(6, 10).Value = 0.00268267035201016
$ws.Cells.Item(6, 11).Value = 0.002588457117764729
$ws.Cells.Item(6, 12).Value = 0.002392056897706832
$ws.Cells.Item(6, 13).Value = 0.002382628539213846
$ws.Cells.Item(6, 14).Value = 0.002311428332587562
$ws.Cells.Item(6, 15).Value = 0.002311428332587562
$ws.Cells.Item(6, 16).Value = 0.002311428332587562
$ws.Cells.Item(6, 17).Value = 0.002288003357419134
$ws.Cells.Item(6, 18).Value = 0.002196440477209957
$ws.Cells.Item(6, 19).Value = 0.002196440477209957
$ws.Cells.Item(6, 20).Value = 0.002196440477209957
$ws.Cells.Item(6, 21).Value = 0.002142897139282686
$ws.Cells.Item(6, 22).Value = 0.002099456016182378
$ws.Cells.Item(6, 23).Value = 0.002099456016182378
$ws.Cells.Item(6, 24).Value = 0.002084485069870205
$ws.Cells.Item(6, 25).Value = 0.002074342538609472

$ws.Cells.Item(7, 3).Value = 0.4459996223449707
$ws.Cells.Item(7, 5).Value = 111.32206186652
$ws.Cells.Item(7, 6).Value = 0.004030069430481654
$ws.Cells.Item(7, 7).Value = 0.003310248356602565
$ws.Cells.Item(7, 8).Value = 0.003150145755943382
$ws.Cells.Item(7, 9).Value = 0.003074561023184581
$ws.Cells.Item(7, 10).Value = 0.002882612388612686
$ws.Cells.Item(7, 11).Value = 0.002717489937048134
$ws.Cells.Item(7, 12).Value = 0.002469547004782736
$ws.Cells.Item(7, 13).Value = 0.002469547004782736
$ws.Cells.Item(7, 14).Value = 0.002469547004782736
$ws.Cells.Item(7, 15).Value = 0.002447003739687123
$ws.Cells.Item(7, 16).Value = 0.002447003739687123
$ws.Cells.Item(7, 17).Value = 0.002435714049703983
$ws.Cells.Item(7, 18).Value = 0.002300014608843351
$ws.Cells.Item(7, 19).Value = 0.002300014608843351
$ws.Cells.Item(7, 20).Value = 0.002300014608843351
$ws.Cells.Item(7, 21).Value = 0.002268460366018334
$ws.Cells.Item(7, 22).Value = 0.002211665355280877
$ws.Cells.Item(7, 23).Value = 0.002210010818828159
$ws.Cells.Item(7, 24).Value = 0.002181766782191171
$ws.Cells.Item(7, 25).Value = 0.002170020699152436

$ws.Cells.Item(8, 3).Value = 0.4430015087127686
$ws.Cells.Item(8, 5).Value = 104.6899488172148
$ws.Cells.Item(8, 6).Value = 0.003834725990376678
$ws.Cells.Item(8, 7).Value = 0.003289797778444041
$ws.Cells.Item(8, 8).Value = 0.003110492103311776
$ws.Cells.Item(8, 9).Value = 0.002875454817562622
$ws.Cells.Item(8, 10).Value = 0.002687072622401896
$ws.Cells.Item(8, 11).Value = 0.002498126010236056
$ws.Cells.Item(8, 12).Value = 0.002485042566229815
$ws.Cells.Item(8, 13).Value = 0.002485042566229815
$ws.Cells.Item(8, 14).Value = 0.002378059387569148
$ws.Cells.Item(8, 15).Value = 0.002378059387569148
$ws.Cells.Item(8, 16).Value = 0.002329469927616456
$ws.Cells.Item(8, 17).Value = 0.00225166085409277
$ws.Cells.Item(8, 18).Value = 0.002238719213278654
$ws.Cells.Item(8, 19).Value = 0.002226332401948779
$ws.Cells.Item(8, 20).Value = 0.002206092285229607
$ws.Cells.Item(8, 21).Value = 0.002196445557850747
$ws.Cells.Item(8, 22).Value = 0.002133707322395599
$ws.Cells.Item(8, 23).Value = 0.002098147596273823
$ws.Cells.Item(8, 24).Value = 0.002074625362361155
$ws.Cells.Item(8, 25).Value = 0.002040739743025629

$ws.Cells.Item(9, 3).Value = 0.4909989833831787
$ws.Cells.Item(9, 5).Value = 106.1216559800596
$ws.Cells.Item(9, 6).Value = 0.003889254253717453
$ws.Cells.Item(9, 7).Value = 0.003014551114137519
$ws.Cells.Item(9, 8).Value = 0.002616799115975948
$ws.Cells.Item(9, 9).Value = 0.002616799115975948
$ws.Cells.Item(9, 10).Value = 0.002599655513651002
$ws.Cells.Item(9, 11).Value = 0.002599655513651002
$ws.Cells.Item(9, 12).Value = 0.002467933114314812
$ws.Cells.Item(9, 13).Value = 0.002397970878992622
$ws.Cells.Item(9, 14).Value = 0.002397970878992622
$ws.Cells.Item(9, 15).Value = 0.002397970878992622
$ws.Cells.Item(9, 16).Value = 0.0023940841595003
$ws.Cells.Item(9, 17).Value = 0.002307245921766953
$ws.Cells.Item(9, 18).Value = 0.002245185974784371
$ws.Cells.Item(9, 19).Value = 0.002245185974784371
$ws.Cells.Item(9, 20).Value = 0.002245185974784371
$ws.Cells.Item(9, 21).Value = 0.002197624923420967
$ws.Cells.Item(9, 22).Value = 0.002144008308945271
$ws.Cells.Item(9, 23).Value = 0.002112538484783031
$ws.Cells.Item(9, 24).Value = 0.00209762301144626
$ws.Cells.Item(9, 25).Value = 0.00206864826471851

$ws.Cells.Item(10, 3).Value = 0.3689980506896973
$ws.Cells.Item(10, 5).Value = 111.7733120226076
$ws.Cells.Item(10, 6).Value = 0.00410297909976543
$ws.Cells.Item(10, 7).Value = 0.00350397393517325
$ws.Cells.Item(10, 8).Value = 0.00307224534411933
$ws.Cells.Item(10, 9).Value = 0.003012399099695675
$ws.Cells.Item(10, 10).Value = 0.002933097295829754
$ws.Cells.Item(10, 11).Value = 0.002915881721768896
$ws.Cells.Item(10, 12).Value = 0.002698718913146217
$ws.Cells.Item(10, 13).Value = 0.002574977175271922
$ws.Cells.Item(10, 14).Value = 0.002574977175271922
$ws.Cells.Item(10, 15).Value = 0.002518827394667479
$ws.Cells.Item(10, 16).Value = 0.002446815587190231
$ws.Cells.Item(10, 17).Value = 0.002446815587190231
$ws.Cells.Item(10, 18).Value = 0.002328634095510031
$ws.Cells.Item(10, 19).Value = 0.002328634095510031
$ws.Cells.Item(10, 20).Value = 0.002293680002360338
$ws.Cells.Item(10, 21).Value = 0.002267667064092304
$ws.Cells.Item(10, 22).Value = 0.002240709441077384
$ws.Cells.Item(10, 23).Value = 0.002229161687121137
$ws.Cells.Item(10, 24).Value = 0.002192257025564604
$ws.Cells.Item(10, 25).Value = 0.002178816998491375

$ws.Cells.Item(11, 3).Value = 0.3689644336700439
$ws.Cells.Item(11, 5).Value = 109.5704456930816
$ws.Cells.Item(11, 6).Value = 0.003999839555189039
$ws.Cells.Item(11, 7).Value = 0.003306768277730246
$ws.Cells.Item(11, 8).Value = 0.002935906117089716
$ws.Cells.Item(11, 9).Value = 0.002775500427575294
$ws.Cells.Item(11, 10).Value = 0.002674540227347458
$ws.Cells.Item(11, 11).Value = 0.002674540227347458
$ws.Cells.Item(11, 12).Value = 0.002629500612150466
$ws.Cells.Item(11, 13).Value = 0.002629500612150466
$ws.Cells.Item(11, 14).Value = 0.002525496840010917
$ws.Cells.Item(11, 15).Value = 0.00241550040729174
$ws.Cells.Item(11, 16).Value = 0.002346541923967417
$ws.Cells.Item(11, 17).Value = 0.002346541923967417
$ws.Cells.Item(11, 18).Value = 0.002278841215346062
$ws.Cells.Item(11, 19).Value = 0.002267045670942911
$ws.Cells.Item(11, 20).Value = 0.002199828668032678
$ws.Cells.Item(11, 21).Value = 0.002199828668032678
$ws.Cells.Item(11, 22).Value = 0.002199828668032678
$ws.Cells.Item(11, 23).Value = 0.002174443175264248
$ws.Cells.Item(11, 24).Value = 0.002152845229098304
$ws.Cells.Item(11, 25).Value = 0.002135876134368063
